# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 01:05"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1666253
$ws.Range("C4").Value = 21159
$ws.Range("D4").Value = 445485
$ws.Range("E4").Value = 1122107
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1014
$ws.Range("H4").Value = 98661

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 347398
$ws.Range("C5").Value = 16508
$ws.Range("D5").Value = 135430
$ws.Range("E5").Value = 189955
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 965
$ws.Range("H5").Value = 22013

# --- Row 51: Panama ---
$ws.Range("B51").Value = 10577
$ws.Range("C51").Value = 310
$ws.Range("D51").Value = 6279
$ws.Range("E51").Value = 3999
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 299

# --- Row 53: Chequia ---
$ws.Range("B53").Value = 8890
$ws.Range("C53").Value = 77
$ws.Range("D53").Value = 6044
$ws.Range("E53").Value = 2532
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 314

# --- Row 55: Noruega ---
$ws.Range("B55").Value = 8346
$ws.Range("C55").Value = 14
$ws.Range("D55").Value = 7727
$ws.Range("E55").Value = 384
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 235

# --- Rows 58/59: Nigeria moves ahead of Marruecos in sort order (Nigeria's
# total cases now exceed Marruecos's), so the country names on these two
# rows swap along with their respective updated data. ---
$ws.Range("A58").Value = "Nigeria"
$ws.Range("B58").Value = 7526
$ws.Range("C58").Value = 265
$ws.Range("D58").Value = 2174
$ws.Range("E58").Value = 5131
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 221

$ws.Range("A59").Value = "Marruecos"
$ws.Range("B59").Value = 7406
$ws.Range("C59").Value = 74
$ws.Range("D59").Value = 4638
$ws.Range("E59").Value = 2570
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 198

# --- Row 179: Angola ---
$ws.Range("B179").Value = 61
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 18
$ws.Range("E179").Value = 39
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 1
$ws.Range("H179").Value = 4
